# Update 16-Jul-2021, end of day update.
# Adds petty-cash entries for 15-Jul-2021 (rows 19-24) and 16-Jul-2021
# (rows 25-34) to "Buku KAS HARIAN"-style ledger on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- existing-row edits (amend 15-Jul-2021 entries) ---------------------
$ws.Range("D19").Formula = "=60000+260000"
$ws.Range("C20").Formula = "=223986500+5943000+20400000"

# --- new rows continuing 15-Jul-2021 -------------------------------------
$ws.Range("B22").Value = "SALES - cash/retail"
$ws.Range("C22").Formula = "=5133025+20463975-20400000"

$ws.Range("B23").Value = "SELISIH - lebih"
$ws.Range("C23").Value = 90000

$ws.Range("B24").Value = "SETOR KE BANK"
$ws.Range("D24").Value = 4000000

# --- 16-Jul-2021 entries ---------------------------------------------------
$ws.Range("A25").Value = 44392
$ws.Range("B25").Value = "Wages Expense"
$ws.Range("D25").Formula = "=60000+260000"

$ws.Range("B26").Value = "TRANSFER BCA"
$ws.Range("D26").Formula = "=1144000+5727000+30500000+7540000+849000+140000+11260000+2390000+1200000"

$ws.Range("B27").Value = "A/R"
$ws.Range("C27").Formula = "=5727000+30500000+7540000+16843000"

$ws.Range("B28").Value = "FREIGHT OUT"
$ws.Range("D28").Formula = "=15000"

$ws.Range("B29").Value = "A/P"
$ws.Range("D29").Formula = "=624000"

$ws.Range("B30").Value = "SALES - cash/retail"
$ws.Range("C30").Formula = "=16718975+2222025-16843000"

$ws.Range("B31").Value = "SELISIH - kurang"
$ws.Range("D31").Value = 99500

$ws.Range("B32").Value = "SETOR KE BANK"
$ws.Range("D32").Value = 1000000

# --- 17-Jul-2021 entry (next day, first line) ------------------------------
$ws.Range("A33").Value = 44393
$ws.Range("B33").Value = "Wages Expense"
$ws.Range("D33").Formula = "=60000"

$ws.Range("B34").Value = "BELI lakban"
$ws.Range("D34").Formula = "=78000"

# --- move the saved selection/scroll to mirror the end of the edit --------
$ws.Range("C53").Select()
